$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.88
$wsSummary.Range("B4").Value = -0.12
$wsSummary.Range("B5").Value = -0.08
$wsSummary.Range("B6").Value = 29
$wsSummary.Range("B7").Value = 9
$wsSummary.Range("B9").Value = 31.03

# --- Sheet: Strategy Status ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.88
$wsStatus.Range("D4").Value = 29
$wsStatus.Range("E4").Value = -0.12
$wsStatus.Range("F4").Value = -0.12
$wsStatus.Range("G4").Value = 31.03

# --- New trade row data (row 30) appended to both "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws) {
    # Column B holds a date-looking string ("2026-02-17") that must stay
    # plain text (matching the rest of the column) instead of being
    # auto-converted to a date serial number by Excel.
    $ws.Range("B30").NumberFormat = "@"
    $ws.Range("B30").Value = "2026-02-17"

    $ws.Range("A30").Value = 29
    $ws.Range("C30").Value = "15:22:15"
    $ws.Range("D30").Value = "MarketMaking"
    $ws.Range("E30").Value = "UP"
    $ws.Range("F30").Value = 0.64
    $ws.Range("G30").Value = 0.8
    $ws.Range("H30").Value = "CLOSED"
    $ws.Range("I30").Value = 25
    $ws.Range("J30").Value = 0.16
    $ws.Range("K30").Value = 99.88
    $ws.Range("L30").Value = 0
    $ws.Range("M30").Value = 0
    $ws.Range("N30").Value = 0.6
    $ws.Range("O30").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P30").Value = "early_exit"
    $ws.Range("Q30").Value = 0.14
}

# --- Sheet: All Trades ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

# --- Sheet: MarketMaking ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
